$wb = $excel.ActiveWorkbook

# =========================================================
# Part 1: update "总计" (summary) sheet - insert 2022-Q3 row
# =========================================================
$wsTotal = $wb.Worksheets.Item(1)

# Shift existing data rows 2..6 down to 3..7 (bottom-up to avoid clobbering)
for ($r = 6; $r -ge 2; $r--) {
    $bVal = $wsTotal.Range("B" + $r).Value2
    $cVal = $wsTotal.Range("C" + $r).Value2
    $dVal = $wsTotal.Range("D" + $r).Value2
    $wsTotal.Range("A" + ($r + 1)).Value = ($r - 1)
    $wsTotal.Range("B" + ($r + 1)).Value = $bVal
    $wsTotal.Range("C" + ($r + 1)).Value = $cVal
    $wsTotal.Range("D" + ($r + 1)).Value = $dVal
}

# Row 7 (A7) is a brand-new cell with no inherited style; copy the style from A6
$wsTotal.Range("A6").Copy()
$wsTotal.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write the new 2022-Q3 summary row into row 2
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.17

# =========================================================
# Part 2: insert a new "2022-Q3" worksheet (before "2022-Q2")
# with its own per-fund holdings detail, cloned formatting
# from the "2022-Q2" sheet
# =========================================================
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# --- row 2 ---
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").NumberFormat = "@"
$wsQ3.Range("B2").Value = '004497'
$wsQ3.Range("B2").Style = "Normal"
$wsQ3.Range("C2").Value = '前海开源多元策略灵活配置混合C'
$wsQ3.Range("D2").NumberFormat = "@"
$wsQ3.Range("D2").Value = '1.68'
$wsQ3.Range("D2").Style = "Normal"
$wsQ3.Range("E2").NumberFormat = "@"
$wsQ3.Range("E2").Value = '93.04'
$wsQ3.Range("E2").Style = "Normal"
$wsQ3.Range("F2").NumberFormat = "@"
$wsQ3.Range("F2").Value = '3.57'
$wsQ3.Range("F2").Style = "Normal"
$wsQ3.Range("G2").NumberFormat = "@"
$wsQ3.Range("G2").Value = '0.0600'
$wsQ3.Range("G2").Style = "Normal"
$wsQ3.Range("H2").Value = 10

# --- row 3 ---
$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").NumberFormat = "@"
$wsQ3.Range("B3").Value = '004496'
$wsQ3.Range("B3").Style = "Normal"
$wsQ3.Range("C3").Value = '前海开源多元策略灵活配置混合A'
$wsQ3.Range("D3").NumberFormat = "@"
$wsQ3.Range("D3").Value = '0.91'
$wsQ3.Range("D3").Style = "Normal"
$wsQ3.Range("E3").NumberFormat = "@"
$wsQ3.Range("E3").Value = '93.04'
$wsQ3.Range("E3").Style = "Normal"
$wsQ3.Range("F3").NumberFormat = "@"
$wsQ3.Range("F3").Value = '3.57'
$wsQ3.Range("F3").Style = "Normal"
$wsQ3.Range("G3").NumberFormat = "@"
$wsQ3.Range("G3").Value = '0.0325'
$wsQ3.Range("G3").Style = "Normal"
$wsQ3.Range("H3").Value = 10

# --- row 4 ---
$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").NumberFormat = "@"
$wsQ3.Range("B4").Value = '011997'
$wsQ3.Range("B4").Style = "Normal"
$wsQ3.Range("C4").Value = '景顺长城安盈回报一年持有期混合A'
$wsQ3.Range("D4").NumberFormat = "@"
$wsQ3.Range("D4").Value = '1.69'
$wsQ3.Range("D4").Style = "Normal"
$wsQ3.Range("E4").NumberFormat = "@"
$wsQ3.Range("E4").Value = '26.07'
$wsQ3.Range("E4").Style = "Normal"
$wsQ3.Range("F4").NumberFormat = "@"
$wsQ3.Range("F4").Value = '1.65'
$wsQ3.Range("F4").Style = "Normal"
$wsQ3.Range("G4").NumberFormat = "@"
$wsQ3.Range("G4").Value = '0.0279'
$wsQ3.Range("G4").Style = "Normal"
$wsQ3.Range("H4").Value = 6

# --- row 5 ---
$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").NumberFormat = "@"
$wsQ3.Range("B5").Value = '014768'
$wsQ3.Range("B5").Style = "Normal"
$wsQ3.Range("C5").Value = '景顺华城稳健6月持有混合C'
$wsQ3.Range("D5").NumberFormat = "@"
$wsQ3.Range("D5").Value = '2.25'
$wsQ3.Range("D5").Style = "Normal"
$wsQ3.Range("E5").NumberFormat = "@"
$wsQ3.Range("E5").Value = '24.78'
$wsQ3.Range("E5").Style = "Normal"
$wsQ3.Range("F5").NumberFormat = "@"
$wsQ3.Range("F5").Value = '1.12'
$wsQ3.Range("F5").Style = "Normal"
$wsQ3.Range("G5").NumberFormat = "@"
$wsQ3.Range("G5").Value = '0.0252'
$wsQ3.Range("G5").Style = "Normal"
$wsQ3.Range("H5").Value = 10

# --- row 6 ---
$wsQ3.Range("A6").Value = 4
$wsQ3.Range("B6").NumberFormat = "@"
$wsQ3.Range("B6").Value = '014767'
$wsQ3.Range("B6").Style = "Normal"
$wsQ3.Range("C6").Value = '景顺华城稳健6月持有混合A'
$wsQ3.Range("D6").NumberFormat = "@"
$wsQ3.Range("D6").Value = '1.26'
$wsQ3.Range("D6").Style = "Normal"
$wsQ3.Range("E6").NumberFormat = "@"
$wsQ3.Range("E6").Value = '24.78'
$wsQ3.Range("E6").Style = "Normal"
$wsQ3.Range("F6").NumberFormat = "@"
$wsQ3.Range("F6").Value = '1.12'
$wsQ3.Range("F6").Style = "Normal"
$wsQ3.Range("G6").NumberFormat = "@"
$wsQ3.Range("G6").Value = '0.0141'
$wsQ3.Range("G6").Style = "Normal"
$wsQ3.Range("H6").Value = 10

# --- row 7 ---
$wsQ3.Range("A7").Value = 5
$wsQ3.Range("B7").NumberFormat = "@"
$wsQ3.Range("B7").Value = '012315'
$wsQ3.Range("B7").Style = "Normal"
$wsQ3.Range("C7").Value = '创金合信港股通成长股票A'
$wsQ3.Range("D7").NumberFormat = "@"
$wsQ3.Range("D7").Value = '0.08'
$wsQ3.Range("D7").Style = "Normal"
$wsQ3.Range("E7").NumberFormat = "@"
$wsQ3.Range("E7").Value = '80.48'
$wsQ3.Range("E7").Style = "Normal"
$wsQ3.Range("F7").NumberFormat = "@"
$wsQ3.Range("F7").Value = '5.00'
$wsQ3.Range("F7").Style = "Normal"
$wsQ3.Range("G7").NumberFormat = "@"
$wsQ3.Range("G7").Value = '0.0040'
$wsQ3.Range("G7").Style = "Normal"
$wsQ3.Range("H7").Value = 7

# --- row 8 ---
$wsQ3.Range("A8").Value = 6
$wsQ3.Range("B8").NumberFormat = "@"
$wsQ3.Range("B8").Value = '012316'
$wsQ3.Range("B8").Style = "Normal"
$wsQ3.Range("C8").Value = '创金合信港股通成长股票C'
$wsQ3.Range("D8").NumberFormat = "@"
$wsQ3.Range("D8").Value = '0.07'
$wsQ3.Range("D8").Style = "Normal"
$wsQ3.Range("E8").NumberFormat = "@"
$wsQ3.Range("E8").Value = '80.48'
$wsQ3.Range("E8").Style = "Normal"
$wsQ3.Range("F8").NumberFormat = "@"
$wsQ3.Range("F8").Value = '5.00'
$wsQ3.Range("F8").Style = "Normal"
$wsQ3.Range("G8").NumberFormat = "@"
$wsQ3.Range("G8").Value = '0.0035'
$wsQ3.Range("G8").Style = "Normal"
$wsQ3.Range("H8").Value = 7

# --- row 9 ---
$wsQ3.Range("A9").Value = 7
$wsQ3.Range("B9").NumberFormat = "@"
$wsQ3.Range("B9").Value = '011998'
$wsQ3.Range("B9").Style = "Normal"
$wsQ3.Range("C9").Value = '景顺长城安盈回报一年持有期混合C'
$wsQ3.Range("D9").NumberFormat = "@"
$wsQ3.Range("D9").Value = '0.08'
$wsQ3.Range("D9").Style = "Normal"
$wsQ3.Range("E9").NumberFormat = "@"
$wsQ3.Range("E9").Value = '26.07'
$wsQ3.Range("E9").Style = "Normal"
$wsQ3.Range("F9").NumberFormat = "@"
$wsQ3.Range("F9").Value = '1.65'
$wsQ3.Range("F9").Style = "Normal"
$wsQ3.Range("G9").NumberFormat = "@"
$wsQ3.Range("G9").Value = '0.0013'
$wsQ3.Range("G9").Style = "Normal"
$wsQ3.Range("H9").Value = 6

# Rows 6-9 are brand-new cells in column A (the source sheet only had rows 2-5);
# copy column-A style from row 2 so they match the others
$wsQ3.Range("A2").Copy()
$wsQ3.Range("A6:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "edit complete"